$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows 2-9: columns E..T updated following Dr Hou's advice
# (ligand/receptor-expressing cell counts 1 -> 3, recomputed downstream values)

$rowData = @{
    2 = @{ E=3; G=36.78646733333333; H=110.359402; I=0.1576941973553631; J=0.1576941973553631;
           K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636;
           Q=238.7643197152762; R=2148.878877437486; S=0.1264970295590645; T=0.1264970295590645 }
    3 = @{ E=3; G=36.78646733333333; H=110.359402; I=0.1576941973553631; J=0.1576941973553631;
           K=3; M=1.600723; N=4.802169; O=0.1978333275383364; P=0.1978333275383364;
           Q=58.88494434921534; R=529.9644991429381; S=0.03119716779629861; T=0.03119716779629862 }
    4 = @{ E=3; G=103.9192913333333; H=311.757874; I=0.445475481188675; J=0.4454754811886751;
           K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636;
           Q=674.4931138852203; R=6070.438024966982; S=0.3573455844083779; T=0.3573455844083779 }
    5 = @{ E=3; G=103.9192913333333; H=311.757874; I=0.445475481188675; J=0.4454754811886751;
           K=3; M=1.600723; N=4.802169; O=0.1978333275383364; P=0.1978333275383364;
           Q=166.3459997809673; R=1497.113998028706; S=0.08812989678029716; T=0.08812989678029716 }
    6 = @{ E=3; G=43.05432033333333; H=129.162961; I=0.1845628835768525; J=0.1845628835768525;
           K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636;
           Q=279.4461183794359; R=2515.015065414923; S=0.1480501941787732; T=0.1480501941787732 }
    7 = @{ E=3; G=43.05432033333333; H=129.162961; I=0.1845628835768525; J=0.1845628835768525;
           K=3; M=1.600723; N=4.802169; O=0.1978333275383364; P=0.1978333275383364;
           Q=68.91804080693433; R=620.262367262409; S=0.03651268939807931; T=0.03651268939807931 }
    8 = @{ E=3; G=49.51716233333334; H=148.551487; I=0.2122674378791094; J=0.2122674378791094;
           K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636;
           Q=321.3935024425713; R=2892.541521983141; S=0.1702738643154481; T=0.1702738643154481 }
    9 = @{ E=3; G=49.51716233333334; H=148.551487; I=0.2122674378791094; J=0.2122674378791094;
           K=3; M=1.600723; N=4.802169; O=0.1978333275383364; P=0.1978333275383364;
           Q=79.26326064170034; R=713.369345775303; S=0.04199357356366132; T=0.04199357356366132 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}
